# ASEAN_HDI.xlsx — add an "additional_info" sheet with HDI classification
# thresholds + a per-country lookup table, matching the commit
# "More Info regarding COVID / Adding more information regarding to covid
# on the final report".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the original sheet Sheet1 -> HDI, then add the new sheet
#    "additional_info" right after it.
# ---------------------------------------------------------------------
$hdi = $wb.Worksheets.Item(1)
$hdi.Name = "HDI"

$info = $wb.Worksheets.Add($null, $hdi)
$info.Name = "additional_info"

# Old selection on the HDI sheet moves from E3 to L7 (no data changed).
$hdi.Range("L7").Select()

# ---------------------------------------------------------------------
# 2. Populate "additional_info".
#    Columns A:B -> HDI classification thresholds (+ "Additional
#    Information" / "Point" mini block further down).
#    Columns D:F -> per-country point + classification lookup (array
#    formula referencing the A:B threshold table).
# ---------------------------------------------------------------------

# NOTE: cell-write order below intentionally matches the original
# authoring sequence so that the regenerated shared-string table comes
# out in the same index order as the source workbook (0-11 already
# exist; 12-24 are newly interned in this exact order: Classification,
# point, classification, Very high/High/Medium/Low, East Asia and the
# Pacific, Developing Country, Brunei, Laos, Additional Information,
# Point).

# --- Header row ---
$info.Range("A1").Value = "Classification"
$info.Range("B1").Value = "point"
$info.Range("D1").Value = "Country"
$info.Range("E1").Value = "point"
$info.Range("F1").Value = "classification"

# --- Classification threshold table (A2:B5) ---
$info.Range("A2").Value = "Very high "
$info.Range("B2").Value = 0.896
$info.Range("A3").Value = "High "
$info.Range("B3").Value = 0.754
$info.Range("A4").Value = "Medium "
$info.Range("B4").Value = 0.636
$info.Range("A5").Value = "Low "
$info.Range("B5").Value = 0.518

# --- Small "Additional Information" block, labels entered first ---
$info.Range("A8").Value = "East Asia and the Pacific"
$info.Range("B8").Value = 0.749
$info.Range("A9").Value = "Developing Country"
$info.Range("B9").Value = 0.685

# --- Per-country point + classification lookup (D2:F11) ---
$countries = @(
    @{ Row = 2;  Name = "Singapore"; Point = 0.939 },
    @{ Row = 3;  Name = "Brunei";    Point = 0.829 },
    @{ Row = 4;  Name = "Malaysia";  Point = 0.803 },
    @{ Row = 5;  Name = "Thailand";  Point = 0.8 },
    @{ Row = 6;  Name = "Indonesia"; Point = 0.705 },
    @{ Row = 7;  Name = "Viet Nam";  Point = 0.703 },
    @{ Row = 8;  Name = "Philippines"; Point = 0.699 },
    @{ Row = 9;  Name = "Laos";      Point = 0.607 },
    @{ Row = 10; Name = "Cambodia";  Point = 0.593 },
    @{ Row = 11; Name = "Myanmar";   Point = 0.585 }
)

foreach ($c in $countries) {
    $r = $c.Row
    $info.Range("D$r").Value = $c.Name
    $info.Range("E$r").Value = $c.Point
    $info.Range("F$r").FormulaArray = "=IF(E$r<B`$5,A`$5,IF(E$r<B`$3,A`$4,IF(E$r<B`$2,A`$3,IF(E$r>B`$2,A`$2,error))))"
}

# --- "Additional Information" header, entered last (matches source
#     shared-string ordering: indices 23 "Additional Information" and
#     24 "Point" come after Brunei/Laos were interned) ---
$info.Range("A7").Value = "Additional Information"
$info.Range("B7").Value = "Point"

# ---------------------------------------------------------------------
# 3. Formatting.
# ---------------------------------------------------------------------

# Thin border + centered header cells (A1:B1)
$hdr = $info.Range("A1:B1")
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108   # xlCenter

# D1:F1 header band — thin border, left aligned
$info.Range("D1:F1").Borders.LineStyle = 1
$info.Range("D1").HorizontalAlignment = -4131   # xlLeft
$info.Range("D1").VerticalAlignment = -4160     # xlTop
$info.Range("E1").HorizontalAlignment = -4131
$info.Range("F1").HorizontalAlignment = -4131
$info.Range("G1").Borders.LineStyle = 1

# Threshold table A2:B5 — thin border, left/top aligned
$info.Range("A2:B5").Borders.LineStyle = 1
$info.Range("A2:A5").HorizontalAlignment = -4131
$info.Range("A2:A5").VerticalAlignment = -4160

# D:F data rows — thin border; E column numeric 0.000 format
$info.Range("D2:F11").Borders.LineStyle = 1
$info.Range("E2:E11").NumberFormat = "0.000"

# Additional-info block A7:B9
$info.Range("A7:B9").Borders.LineStyle = 1
$info.Range("A7").Font.Name = "Arial"
$info.Range("A7").WrapText = $true
$info.Range("A7").HorizontalAlignment = -4131
$info.Range("A7").VerticalAlignment = -4108     # xlCenter

$info.Range("A8:A9").HorizontalAlignment = -4108
$info.Range("A8:A9").VerticalAlignment = -4160

$info.Range("B8:B9").NumberFormat = "#,###,##0.000"
$info.Range("B8:B9").Font.Name = "Arial"
$info.Range("B8:B9").HorizontalAlignment = -4108
$info.Range("B8:B9").VerticalAlignment = -4108

# ---------------------------------------------------------------------
# 4. Column widths (best-fit approximations of the original workbook).
# ---------------------------------------------------------------------
$info.Columns.Item(1).ColumnWidth = 29.11
$info.Columns.Item(3).ColumnWidth = 11.55
$info.Columns.Item(4).ColumnWidth = 29.11
$info.Columns.Item(5).ColumnWidth = 5.66
$info.Columns.Item(6).ColumnWidth = 11.55
$info.Columns.Item(7).ColumnWidth = 17.33

# ---------------------------------------------------------------------
# 5. Sheet view / selection / zoom, then activate "additional_info".
# ---------------------------------------------------------------------
$info.PageSetup.Orientation = 1   # xlPortrait

$info.Range("C8").Select()
$excel.ActiveWindow.Zoom = 145

$info.Select()
